$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows below the current CruizcoreGyro row (row 15) so that the
# sensor now occupies three rows (15-17), mirroring the multi-mode layout used
# by other multi-mode sensors (e.g. AccelHTSensor rows 2-4).
$ws.Rows.Item(16).Resize(2).Insert()

# Row 15: CruizcoreGyro - add the "new name", first mode (Acceleration) and the
# (non-hyperlinked) spec URL.
$ws.Range("B15").Value = "CruizcoreGyro"
$ws.Range("F15").Value = "Acceleration"
$ws.Range("H15").Value = "http://xgl.minfinity.com/Downloads/Downloads.html"

# Row 16: second mode (Rate) for CruizcoreGyro.
$ws.Range("F16").Value = "Rate"
$ws.Range("G16").Value = "SampleProvider"

# Row 17: third mode (Angle) for CruizcoreGyro.
$ws.Range("F17").Value = "Angle"
$ws.Range("G17").Value = "SampleProvider"

# Row 24 (was row 22 before the insert above): DThermalIR gains a claimed
# "new name", developer, finished flag, mode and interface.
$ws.Range("B24").Value = "DexterThermalIR"
$ws.Range("D24").Value = "Lawrie"
$ws.Range("E24").Value = "N"
$ws.Range("F24").Value = "Temperature"
$ws.Range("G24").Value = "SampleProvider"

# Restore the active selection to match the author's final cursor position.
$ws.Range("D20").Select()
